# Fruta / hortaliza, semanal
# Insert a new data row before row 87 (pushes existing rows 87-145 down to 88-146)
# and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87 (shifts rows 87..145 down to 88..146).
$ws.Rows(87).Insert()

$ws.Range("A87").Value = 7
$ws.Range("B87").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C87").Value = "Ñuble"
$ws.Range("D87").Value = "2022-04-22"
$ws.Range("E87").Value = 16
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100109
$ws.Range("H87").Value = "Uva"
$ws.Range("I87").Value = 100109001
$ws.Range("J87").Value = "Uva"
$ws.Range("K87").Value = "Red Globe"
$ws.Range("L87").Value = "Primera"
$ws.Range("M87").Value = 120
$ws.Range("N87").Value = 8000
$ws.Range("O87").Value = 9000
$ws.Range("P87").Value = 8500
$ws.Range("Q87").Value = "`$/bandeja 18 kilos"
$ws.Range("R87").Value = "Región de O'Higgins"
$ws.Range("S87").Value = 472
$ws.Range("T87").Value = 18
